$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1) Status ("A") column corrections on a handful of existing rows
#    (these precede the two brand-new TODO rows and are independent
#    single-cell edits).
# ------------------------------------------------------------------
$ws.Range("A62").Value = "Done"          # Molten salt linear fresnel - match with script from Ty
$ws.Range("A89").Value = "Future"        # CEC Inverter updates
$ws.Range("A90").Value = "Done"          # PBNS update for dispatch factors

# ------------------------------------------------------------------
# 2) Insert two brand-new TODO rows right after row 90 (pushes the
#    "Pre-release" CSP-constants rows, and everything below, down by
#    two rows).
# ------------------------------------------------------------------
$ws.Rows("91:92").Insert()

$ws.Range("A91").Value = "Not done"
$ws.Range("B91").Value = "Move all CSP constants from ui to respective compute modules"
$ws.Range("C91").Value = "Steve"

$ws.Range("A92").Value = "Not done"
$ws.Range("B92").Value = "Fix issue with parametric grid editor read only display of monthly and hourly outputs."

# ------------------------------------------------------------------
# 3) Conditional-formatting ranges shift down by two rows to keep
#    pointing at the same logical blocks of data.
# ------------------------------------------------------------------
$ws.Range("A1:A93 A96:A1048576").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("A1:A95 A98:A1048576"))
$ws.Range("A1:A93 A96:A1048576").FormatConditions.Item(2).ModifyAppliesToRange($ws.Range("A1:A95 A98:A1048576"))
$ws.Range("A1:A93 A96:A1048576").FormatConditions.Item(3).ModifyAppliesToRange($ws.Range("A1:A95 A98:A1048576"))
$ws.Range("A1:A93 A96:A1048576").FormatConditions.Item(4).ModifyAppliesToRange($ws.Range("A1:A95 A98:A1048576"))
$ws.Range("A1:A93 A96:A1048576").FormatConditions.Item(5).ModifyAppliesToRange($ws.Range("A1:A95 A98:A1048576"))

$ws.Range("A94:A95").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("A96:A97"))
$ws.Range("A94:A95").FormatConditions.Item(2).ModifyAppliesToRange($ws.Range("A96:A97"))
$ws.Range("A94:A95").FormatConditions.Item(3).ModifyAppliesToRange($ws.Range("A96:A97"))
$ws.Range("A94:A95").FormatConditions.Item(4).ModifyAppliesToRange($ws.Range("A96:A97"))
$ws.Range("A94:A95").FormatConditions.Item(5).ModifyAppliesToRange($ws.Range("A96:A97"))

# ------------------------------------------------------------------
# 4) Restore the view: frozen pane scrolled to A59, selection on B92.
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("B92").Select()
$excel.ActiveWindow.ScrollRow = 59
